$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 14 from 45203 to 45204
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45204
    }
}
